# Update "想去人数" (number of people interested) counts across the
# "展览" (Exhibitions), "演出" (Performances) and "全部类型" (All types)
# worksheets to reflect the latest generated numbers.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1651
$ws1.Range("F3").Value  = 9062
$ws1.Range("F5").Value  = 505
$ws1.Range("F6").Value  = 695
$ws1.Range("F7").Value  = 865
$ws1.Range("F8").Value  = 190
$ws1.Range("F9").Value  = 52
$ws1.Range("F10").Value = 86
$ws1.Range("F11").Value = 5674
$ws1.Range("F13").Value = 382
$ws1.Range("F15").Value = 4357
$ws1.Range("F19").Value = 17
$ws1.Range("F22").Value = 252
$ws1.Range("F24").Value = 2703

# --- Sheet: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 10

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1651
$ws4.Range("F3").Value  = 9062
$ws4.Range("F5").Value  = 10
$ws4.Range("F6").Value  = 505
$ws4.Range("F7").Value  = 695
$ws4.Range("F8").Value  = 865
$ws4.Range("F9").Value  = 190
$ws4.Range("F10").Value = 52
$ws4.Range("F11").Value = 86
$ws4.Range("F12").Value = 5674
$ws4.Range("F14").Value = 382
$ws4.Range("F16").Value = 4357
$ws4.Range("F20").Value = 17
$ws4.Range("F23").Value = 252
$ws4.Range("F25").Value = 2703
